$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.503.27"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "'2.378.18"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'506.18"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'130.76"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").Value = "'2.390.16"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").Value = "'0.0986"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  +5.97%  "
$ws.Range("D13").Value = "'0.330"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").Value = "'2.800.69"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "'56.584.23"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "'21.65"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "'2.383.39"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "'309.66"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "'6.29"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'66.37"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "'0.996"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").Value = "'173.84"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").Value = "'5.85"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("D36").Value = "'17.68"
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("D37").Value = "'1.19"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  -3.07%  "
$ws.Range("D39").Value = "'0.825"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'127.61"
$ws.Range("E43").Value = "  -4.16%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'4.76"
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "'0.0898"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "'240.75"
$ws.Range("E47").Value = "  -4.57%  "
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").Value = "'17.08"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  -0.19%  "
